$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates (column E "Information Extraction Method", plus C5/D5) ---
$ws.Range("E3").Value = "Direct mapping / 1000"
$ws.Range("E4").Value = "preprocessing table"

$ws.Range("C5").Value = "PV1-9.7"
$ws.Range("D5").Value = "Consulting Doctor Degree"
$ws.Range("E5").Value = "preprocessing table (hl7 docs)"

$ws.Range("E6").Value = "preprocessing table (hl7 docs)"
$ws.Range("E8").Value = "preprocessing table"
$ws.Range("E9").Value = "preprocessing table (hl7 docs)"
$ws.Range("E10").Value = "preprocessing table"
$ws.Range("E12").Value = "Direct mapping * 86400"

# --- New column I: width + a formatted (empty) cell at I2 ---
$ws.Columns.Item(9).ColumnWidth = 12.14
$ws.Range("I2").HorizontalAlignment = -4131
$ws.Range("I2").VerticalAlignment = -4108
$ws.Range("I2").WrapText = $true

# --- Selection moved to E7 ---
$ws.Range("E7").Select()
